$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 42.517849
$ws.Range("H2").Value = 127.553547
$ws.Range("I2").Value = 0.02311569285614191
$ws.Range("J2").Value = 0.02311569285614191
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 25.94532566666667
$ws.Range("N2").Value = 77.835977
$ws.Range("O2").Value = 0.5401813355606462
$ws.Range("P2").Value = 0.5401813355606462
$ws.Range("Q2").Value = 1103.139438951158
$ws.Range("R2").Value = 9928.254950560418
$ws.Range("S2").Value = 0.01248666583944043
$ws.Range("T2").Value = 0.01248666583944043

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 42.517849
$ws.Range("H3").Value = 127.553547
$ws.Range("I3").Value = 0.02311569285614191
$ws.Range("J3").Value = 0.02311569285614191
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 18.51427066666667
$ws.Range("N3").Value = 55.542812
$ws.Range("O3").Value = 0.3854668692210787
$ws.Range("P3").Value = 0.3854668692210786
$ws.Range("Q3").Value = 787.1869645504627
$ws.Range("R3").Value = 7084.682680954164
$ws.Range("S3").Value = 0.008910333755133079
$ws.Range("T3").Value = 0.008910333755133075

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 42.517849
$ws.Range("H4").Value = 127.553547
$ws.Range("I4").Value = 0.02311569285614191
$ws.Range("J4").Value = 0.02311569285614191
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.571174000000001
$ws.Range("N4").Value = 10.713522
$ws.Range("O4").Value = 0.07435179521827505
$ws.Range("P4").Value = 0.07435179521827504
$ws.Range("Q4").Value = 151.838636884726
$ws.Range("R4").Value = 1366.547731962534
$ws.Range("S4").Value = 0.001718693261568407
$ws.Range("T4").Value = 0.001718693261568407

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1644.738728666666
$ws.Range("H5").Value = 4934.216186
$ws.Range("I5").Value = 0.8941956419399297
$ws.Range("J5").Value = 0.8941956419399296
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 25.94532566666667
$ws.Range("N5").Value = 77.835977
$ws.Range("O5").Value = 0.5401813355606462
$ws.Range("P5").Value = 0.5401813355606462
$ws.Range("Q5").Value = 42673.28195183596
$ws.Range("R5").Value = 384059.5375665237
$ws.Range("S5").Value = 0.4830277961156206
$ws.Range("T5").Value = 0.4830277961156205

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1644.738728666666
$ws.Range("H6").Value = 4934.216186
$ws.Range("I6").Value = 0.8941956419399297
$ws.Range("J6").Value = 0.8941956419399296
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 18.51427066666667
$ws.Range("N6").Value = 55.542812
$ws.Range("O6").Value = 0.3854668692210787
$ws.Range("P6").Value = 0.3854668692210786
$ws.Range("Q6").Value = 30451.13799848389
$ws.Range("R6").Value = 274060.241986355
$ws.Range("S6").Value = 0.3446827945697174
$ws.Range("T6").Value = 0.3446827945697173

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1644.738728666666
$ws.Range("H7").Value = 4934.216186
$ws.Range("I7").Value = 0.8941956419399297
$ws.Range("J7").Value = 0.8941956419399296
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.571174000000001
$ws.Range("N7").Value = 10.713522
$ws.Range("O7").Value = 0.07435179521827505
$ws.Range("P7").Value = 0.07435179521827504
$ws.Range("Q7").Value = 5873.648184607455
$ws.Range("R7").Value = 52862.83366146709
$ws.Range("S7").Value = 0.06648505125459166
$ws.Range("T7").Value = 0.06648505125459163

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 152.093394
$ws.Range("H8").Value = 456.280182
$ws.Range("I8").Value = 0.08268866520392831
$ws.Range("J8").Value = 0.0826886652039283
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 25.94532566666667
$ws.Range("N8").Value = 77.835977
$ws.Range("O8").Value = 0.5401813355606462
$ws.Range("P8").Value = 0.5401813355606462
$ws.Range("Q8").Value = 3946.112639078646
$ws.Range("R8").Value = 35515.01375170782
$ws.Range("S8").Value = 0.04466687360558513
$ws.Range("T8").Value = 0.04466687360558512

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 152.093394
$ws.Range("H9").Value = 456.280182
$ws.Range("I9").Value = 0.08268866520392831
$ws.Range("J9").Value = 0.0826886652039283
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 18.51427066666667
$ws.Range("N9").Value = 55.542812
$ws.Range("O9").Value = 0.3854668692210787
$ws.Range("P9").Value = 0.3854668692210786
$ws.Range("Q9").Value = 2815.898263127976
$ws.Range("R9").Value = 25343.08436815179
$ws.Range("S9").Value = 0.03187374089622819
$ws.Range("T9").Value = 0.03187374089622819

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 152.093394
$ws.Range("H10").Value = 456.280182
$ws.Range("I10").Value = 0.08268866520392831
$ws.Range("J10").Value = 0.0826886652039283
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.571174000000001
$ws.Range("N10").Value = 10.713522
$ws.Range("O10").Value = 0.07435179521827505
$ws.Range("P10").Value = 0.07435179521827504
$ws.Range("Q10").Value = 543.1519742245562
$ws.Range("R10").Value = 4888.367768021005
$ws.Range("S10").Value = 0.006148050702114984
$ws.Range("T10").Value = 0.006148050702114981
